# Auto-generated edit script: updates market-price derived columns (H-N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# refreshed marketboard data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = [ordered]@{
    "H5" = 162.83333
    "I5" = 203.75
    "J5" = 81
    "K5" = 203.75
    "L5" = 81
    "M5" = -88.75
    "N5" = -311
    "H6" = 989
    "I6" = 4
    "J6" = 1580
    "K6" = 12
    "L6" = 4740
    "M6" = 100
    "N6" = -4964
    "H33" = 319.2
    "I33" = 330.33334
    "K33" = 330.33334
    "M33" = -101.33334
    "H64" = 5187.5
    "I64" = 4537.4
    "K64" = 4537.4
    "M64" = -4289.4
    "H67" = 5187.5
    "I67" = 4537.4
    "K67" = 4537.4
    "M67" = -3679.4
    "H76" = 25003050
    "I76" = 33335966
    "J76" = 4300
    "K76" = 33335966
    "L76" = 4300
    "M76" = -33335651
    "N76" = -4930
    "H79" = 25003050
    "I79" = 33335966
    "J79" = 4300
    "K79" = 33335966
    "L79" = 4300
    "M79" = -33334874
    "N79" = -6484
    "H80" = 838.7857
    "I80" = 508.14285
    "K80" = 1524.42855
    "M80" = -526.4285500000001
    "H83" = 838.7857
    "I83" = 508.14285
    "K83" = 4573.28565
    "M83" = 418.7143500000002
    "H86" = 12167
    "I86" = 10748.5
    "K86" = 10748.5
    "M86" = -9625.5
    "H87" = 89998
    "I87" = 0
    "K87" = 0
    "M87" = $null
    "H89" = 12167
    "I89" = 10748.5
    "K89" = 53742.5
    "M89" = -48126.5
    "H90" = 89998
    "I90" = 0
    "K90" = 0
    "M90" = $null
    "H98" = 1020.05554
    "I98" = 622.625
    "J98" = 4199.5
    "K98" = 622.625
    "L98" = 4199.5
    "M98" = 875.375
    "N98" = -7195.5
    "H122" = 1020.05554
    "I122" = 622.625
    "J122" = 4199.5
    "K122" = 1867.875
    "L122" = 12598.5
    "M122" = 582.125
    "N122" = -17498.5
    "H137" = 3583.423
    "I137" = 3219.3333
    "J137" = 7952.5
    "K137" = 9657.999899999999
    "L137" = 23857.5
    "M137" = -7107.999899999999
    "N137" = -28957.5
}
foreach ($key in $ALC_updates.Keys) {
    $ws.Range($key).Value = $ALC_updates[$key]
}

$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = [ordered]@{
    "H2" = 2173.4443
    "J2" = 2302.7778
    "L2" = 2302.7778
    "N2" = -2528.7778
    "H32" = 1837.6207
    "I32" = 1546.1428
    "K32" = 1546.1428
    "M32" = -1259.1428
    "H45" = 2285.4614
    "I45" = 2294.2
    "J45" = 2256.3333
    "K45" = 2294.2
    "L45" = 2256.3333
    "M45" = -1917.2
    "N45" = -3010.3333
    "H63" = 11173.941
    "I63" = 10518.429
    "J63" = 14233
    "K63" = 10518.429
    "L63" = 14233
    "M63" = -9832.429
    "N63" = -15605
    "H66" = 11173.941
    "I66" = 10518.429
    "J66" = 14233
    "K66" = 52592.145
    "L66" = 71165
    "M66" = -49160.145
    "N66" = -78029
    "H94" = 89998
    "J94" = 89998
    "L94" = 89998
    "N94" = -91800
    "H116" = 2173.4443
    "J116" = 2302.7778
    "L116" = 2302.7778
    "N116" = -6890.7778
    "H122" = 4505.143
    "I122" = 3714.4167
    "K122" = 11143.2501
    "M122" = -8693.250100000001
}
foreach ($key in $ARM_updates.Keys) {
    $ws.Range($key).Value = $ARM_updates[$key]
}

$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = [ordered]@{
    "H3" = 2173.4443
    "J3" = 2302.7778
    "L3" = 2302.7778
    "N3" = -2530.7778
    "H86" = 4654.4443
    "I86" = 2572.0667
    "J86" = 15066.333
    "K86" = 2572.0667
    "L86" = 15066.333
    "M86" = -1449.0667
    "N86" = -17312.333
    "H89" = 4654.4443
    "I89" = 2572.0667
    "J89" = 15066.333
    "K89" = 12860.3335
    "L89" = 75331.66500000001
    "M89" = -7244.333499999999
    "N89" = -86563.66500000001
}
foreach ($key in $BSM_updates.Keys) {
    $ws.Range($key).Value = $BSM_updates[$key]
}

$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = [ordered]@{
    "H86" = 9141.714
    "J86" = 8749.75
    "L86" = 8749.75
    "N86" = -10995.75
    "H89" = 9141.714
    "J89" = 8749.75
    "L89" = 43748.75
    "N89" = -54980.75
    "H132" = 1493.6666
    "J132" = 0
    "L132" = 0
    "N132" = $null
}
foreach ($key in $CRP_updates.Keys) {
    $ws.Range($key).Value = $CRP_updates[$key]
}

$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = [ordered]@{
    "H16" = 1850
    "J16" = 1850
    "L16" = 5550
    "N16" = -5896
    "H23" = 362.85715
    "I23" = 277
    "J23" = 397.2
    "K23" = 831
    "L23" = 1191.6
    "M23" = -596
    "N23" = -1661.6
    "H38" = 226.16667
    "I38" = 290
    "J38" = 98.5
    "K38" = 870
    "L38" = 295.5
    "M38" = -523
    "N38" = -989.5
    "H62" = 8396.75
    "J62" = 7662.3335
    "L62" = 22987.0005
    "N62" = -24359.0005
    "H65" = 8396.75
    "J65" = 7662.3335
    "L65" = 68961.0015
    "N65" = -75825.0015
    "H121" = 834271.7
    "J121" = 1251002.9
    "L121" = 3753008.7
    "N121" = -3755628.7
    "H132" = 1909.5454
    "I132" = 1878.3334
    "K132" = 16905.0006
    "M132" = -14375.0006
}
foreach ($key in $CUL_updates.Keys) {
    $ws.Range($key).Value = $CUL_updates[$key]
}

$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = [ordered]@{
    "H2" = 125.6
    "I2" = 142.92308
    "J2" = 13
    "K2" = 142.92308
    "L2" = 13
    "M2" = -29.92308
    "N2" = -239
    "H70" = 8812.666999999999
    "J70" = 8979
    "L70" = 8979
    "N70" = -9519
    "H73" = 8812.666999999999
    "J73" = 8979
    "L73" = 8979
    "N73" = -10851
    "H80" = 3017.25
    "I80" = 2702
    "K80" = 2702
    "M80" = -1704
    "H83" = 3017.25
    "I83" = 2702
    "K83" = 13510
    "M83" = -8518
    "H132" = 5041.9165
    "I132" = 4701.3335
    "K132" = 14104.0005
    "M132" = -11574.0005
    "H134" = 44666.668
    "J134" = 44666.668
    "L134" = 134000.004
    "N134" = -139070.004
}
foreach ($key in $GSM_updates.Keys) {
    $ws.Range($key).Value = $GSM_updates[$key]
}

$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = [ordered]@{
    "H46" = 2055
    "I46" = 1317
    "J46" = 3285
    "K46" = 1317
    "L46" = 3285
    "M46" = -1129
    "N46" = -3661
    "H108" = 48990
    "J108" = 48990
    "L108" = 48990
    "N108" = -56670
    "H136" = 25002924
    "I136" = 3013
    "J136" = 62502790
    "K136" = 9039
    "L136" = 187508370
    "M136" = -6489
    "N136" = -187513470
}
foreach ($key in $LTW_updates.Keys) {
    $ws.Range($key).Value = $LTW_updates[$key]
}

$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = [ordered]@{
    "H136" = 9882.5
    "I136" = 9882.5
    "K136" = 29647.5
    "M136" = -27097.5
}
foreach ($key in $WVR_updates.Keys) {
    $ws.Range($key).Value = $WVR_updates[$key]
}
